$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to track "Depression" scores for 2023-01-05..09; it now
# tracks "Stress" scores for a new week, 2024-05-22..29 (9 rows total,
# up from 6).
$dates  = @("2024-05-22","2024-05-23","2024-05-24","2024-05-25","2024-05-26","2024-05-27","2024-05-28","2024-05-29")
$scores = @(5,8,7,7,6,6,7,5)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("B$row").Value = "Stress"
    $ws.Range("C$row").Value = $scores[$i]
}

# The first score (C2) carries an underline.
$ws.Range("C2").Font.Underline = $true

# Column A now holds dates, so size it like the rest of the bestFit columns.
$ws.Columns("A").ColumnWidth = 10.1640625

# Leave the selection where the user was last positioned (just past the
# new data).
$ws.Range("C10").Select()
